$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update headers for columns E, F, G
$ws.Range("E1").Value = "Max_Absorption_nm"
$ws.Range("F1").Value = "Max_f_osc"
$ws.Range("G1").Value = "Max_Excitation_eV"

# Update data rows: E = old Max_Absorption_nm (col I), F = old Max_f_osc (col J), G = new Max_Excitation_eV
$ws.Range("E2").Value = 585
$ws.Range("F2").Value = 0.9257919999999999
$ws.Range("G2").Value = 2.12
$ws.Range("E3").Value = 576
$ws.Range("F3").Value = 0.898504
$ws.Range("G3").Value = 2.15
$ws.Range("E4").Value = 938
$ws.Range("F4").Value = 0.748518
$ws.Range("G4").Value = 1.32
$ws.Range("E5").Value = 934
$ws.Range("F5").Value = 0.68387
$ws.Range("G5").Value = 1.33
$ws.Range("E6").Value = 918
$ws.Range("F6").Value = 0.616428
$ws.Range("G6").Value = 1.35
$ws.Range("E7").Value = 530
$ws.Range("F7").Value = 0.912547
$ws.Range("G7").Value = 2.34
$ws.Range("E8").Value = 990
$ws.Range("F8").Value = 0.598239
$ws.Range("G8").Value = 1.25
$ws.Range("E9").Value = 642
$ws.Range("F9").Value = 0.6543369999999999
$ws.Range("G9").Value = 1.93
$ws.Range("E10").Value = 683
$ws.Range("F10").Value = 1.467794
$ws.Range("G10").Value = 1.82
$ws.Range("E11").Value = 415
$ws.Range("F11").Value = 0.140013
$ws.Range("G11").Value = 2.99
$ws.Range("E12").Value = 411
$ws.Range("F12").Value = 0.262642
$ws.Range("G12").Value = 3.02
$ws.Range("E13").Value = 358
$ws.Range("F13").Value = 0.109675
$ws.Range("G13").Value = 3.46
$ws.Range("E14").Value = 474
$ws.Range("F14").Value = 0.498719
$ws.Range("G14").Value = 2.61
$ws.Range("E15").Value = 519
$ws.Range("F15").Value = 0.533673
$ws.Range("G15").Value = 2.39
$ws.Range("E16").Value = 342
$ws.Range("F16").Value = 0.321588
$ws.Range("G16").Value = 3.62
$ws.Range("E17").Value = 304
$ws.Range("F17").Value = 0.321389
$ws.Range("G17").Value = 4.08
$ws.Range("E18").Value = 328
$ws.Range("F18").Value = 0.453893
$ws.Range("G18").Value = 3.78
$ws.Range("E19").Value = 517
$ws.Range("F19").Value = 0.55103
$ws.Range("G19").Value = 2.4
$ws.Range("E20").Value = 273
$ws.Range("F20").Value = 0.145043
$ws.Range("G20").Value = 4.55
$ws.Range("E21").Value = 473
$ws.Range("F21").Value = 0.379511
$ws.Range("G21").Value = 2.62
$ws.Range("E22").Value = 516
$ws.Range("F22").Value = 0.948782
$ws.Range("G22").Value = 2.4
$ws.Range("E23").Value = 474
$ws.Range("F23").Value = 0.708416
$ws.Range("G23").Value = 2.61
$ws.Range("E24").Value = 552
$ws.Range("F24").Value = 1.03911
$ws.Range("G24").Value = 2.25
$ws.Range("E25").Value = 584
$ws.Range("F25").Value = 0.695781
$ws.Range("G25").Value = 2.12

# Remove now-unused columns H, I, J (Total_Energy_Hartree, Solvation_Energy_eV, Surface_Area_A2, Molecular_Volume_A3 data no longer present)
$ws.Range("H1:J25").Clear()

Write-Output "edit applied"
